$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet 1: LP1912 ---
$ws1.Range("A2").Value = "Última actualización: 20:47:24"
$ws1.Range("A3").Value = "Total filas: 484"
$ws1.Range("A35").Value = "05:49:10"
$ws1.Range("C35").Value = "23_HERNANDEZ"
$ws1.Range("D35").Value = 76
$ws1.Range("A36").Value = "05:19:24"
$ws1.Range("C36").Value = "15_ABASTO"
$ws1.Range("D36").Value = 106
$ws1.Range("A80").Value = "07:19:37"
$ws1.Range("C80").Value = "14_ABASTO"
$ws1.Range("D80").Value = 84
$ws1.Range("A81").Value = "08:19:33"
$ws1.Range("C81").Value = "16_SANTA ANA"
$ws1.Range("D81").Value = 24
$ws1.Range("A140").Value = "10:52:48"
$ws1.Range("C140").Value = "10_OLMOS"
$ws1.Range("D140").Value = 0
$ws1.Range("A141").Value = "10:12:35"
$ws1.Range("C141").Value = "15_ABASTO"
$ws1.Range("D141").Value = 40
$ws1.Range("C156").Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Range("C157").Value = "16_SANTA ANA"
$ws1.Range("A186").Value = "11:17:08"
$ws1.Range("C186").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D186").Value = 50
$ws1.Range("A187").Value = "11:59:06"
$ws1.Range("C187").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D187").Value = 8
$ws1.Range("C188").Value = "14_ABASTO"
$ws1.Range("C202").Value = "23_HERNANDEZ"
$ws1.Range("C203").Value = "11_ETCHEVERRY"
$ws1.Range("A261").Value = "14:44:25"
$ws1.Range("C261").Value = "23_HERNANDEZ"
$ws1.Range("D261").Value = 20
$ws1.Range("A262").Value = "14:00:52"
$ws1.Range("C262").Value = "10_OLMOS"
$ws1.Range("D262").Value = 64
$ws1.Range("A307").Value = "15:51:48"
$ws1.Range("C307").Value = "225_GOMEZ"
$ws1.Range("D307").Value = 51
$ws1.Range("A308").Value = "14:44:25"
$ws1.Range("C308").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D308").Value = 118
$ws1.Range("C321").Value = "23_HERNANDEZ"
$ws1.Range("C323").Value = "215A_EL PATO"
$ws1.Range("A337").Value = "17:15:09"
$ws1.Range("C337").Value = "10_OLMOS"
$ws1.Range("D337").Value = 19
$ws1.Range("A338").Value = "16:52:47"
$ws1.Range("C338").Value = "23_HERNANDEZ"
$ws1.Range("D338").Value = 42
$ws1.Range("A368").Value = "18:19:32"
$ws1.Range("C368").Value = "16_SANTA ANA"
$ws1.Range("D368").Value = 1
$ws1.Range("A369").Value = "17:57:54"
$ws1.Range("C369").Value = "26_HERNANDEZ"
$ws1.Range("D369").Value = 23
$ws1.Range("A399").Value = "18:37:39"
$ws1.Range("C399").Value = "23_HERNANDEZ"
$ws1.Range("D399").Value = 28
$ws1.Range("A400").Value = "17:15:09"
$ws1.Range("C400").Value = "11_ETCHEVERRY"
$ws1.Range("D400").Value = 110
$ws1.Range("C412").Value = "14_ABASTO"
$ws1.Range("A413").Value = "17:42:01"
$ws1.Range("C413").Value = "26_HERNANDEZ"
$ws1.Range("D413").Value = 99
$ws1.Range("A414").Value = "18:49:07"
$ws1.Range("C414").Value = "16_SANTA ANA"
$ws1.Range("D414").Value = 32
$ws1.Range("A444").Value = "19:39:27"
$ws1.Range("C444").Value = "16_SANTA ANA"
$ws1.Range("D444").Value = 43
$ws1.Range("A445").Value = "18:37:39"
$ws1.Range("C445").Value = "11_ETCHEVERRY"
$ws1.Range("D445").Value = 105
$ws1.Range("A459").Value = "19:15:23"
$ws1.Range("C459").Value = "10_OLMOS"
$ws1.Range("D459").Value = 101
$ws1.Range("A460").Value = "19:39:27"
$ws1.Range("C460").Value = "27_EL RETIRO"
$ws1.Range("D460").Value = 77
$ws1.Range("A462").Value = "20:47:24"
$ws1.Range("B462").Value = "21:01"
$ws1.Range("C462").Value = "16_SANTA ANA"
$ws1.Range("D462").Value = 14
$ws1.Range("A463").Value = "19:15:23"
$ws1.Range("B463").Value = "21:04"
$ws1.Range("C463").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D463").Value = 109
$ws1.Range("A464").Value = "19:52:49"
$ws1.Range("B464").Value = "21:07"
$ws1.Range("D464").Value = 75
$ws1.Range("A465").Value = "19:15:23"
$ws1.Range("B465").Value = "21:08"
$ws1.Range("C465").Value = "215B_EL PATO"
$ws1.Range("D465").Value = 113
$ws1.Range("A466").Value = "20:34:16"
$ws1.Range("B466").Value = "21:16"
$ws1.Range("C466").Value = "16_SANTA ANA"
$ws1.Range("D466").Value = 42
$ws1.Range("A467").Value = "19:52:49"
$ws1.Range("B467").Value = "21:20"
$ws1.Range("D467").Value = 88
$ws1.Range("A468").Value = "19:39:27"
$ws1.Range("B468").Value = "21:21"
$ws1.Range("C468").Value = "26_HERNANDEZ"
$ws1.Range("D468").Value = 102
$ws1.Range("A469").Value = "19:52:49"
$ws1.Range("C469").Value = "10_OLMOS"
$ws1.Range("D469").Value = 90
$ws1.Range("A470").Value = "20:12:55"
$ws1.Range("B470").Value = "21:22"
$ws1.Range("D470").Value = 70
$ws1.Range("A471").Value = "20:34:16"
$ws1.Range("C471").Value = "15_ABASTO"
$ws1.Range("D471").Value = 49
$ws1.Range("A472").Value = "19:39:27"
$ws1.Range("B472").Value = "21:23"
$ws1.Range("C472").Value = "10_OLMOS"
$ws1.Range("D472").Value = 104
$ws1.Range("A473").Value = "20:12:55"
$ws1.Range("C473").Value = "14_ABASTO"
$ws1.Range("D473").Value = 85
$ws1.Range("B474").Value = "21:37"
$ws1.Range("C474").Value = "17_ROMERO"
$ws1.Range("D474").Value = 105
$ws1.Range("A476").Value = "19:52:49"
$ws1.Range("B476").Value = "21:38"
$ws1.Range("C476").Value = "14_ABASTO"
$ws1.Range("D476").Value = 106
$ws1.Range("A477").Value = "20:12:55"
$ws1.Range("B477").Value = "21:46"
$ws1.Range("D477").Value = 94
$ws1.Range("A478").Value = "19:52:49"
$ws1.Range("B478").Value = "21:47"
$ws1.Range("C478").Value = "215A_EL PATO"
$ws1.Range("A479").Value = "20:12:55"
$ws1.Range("B479").Value = "22:07"
$ws1.Range("C479").Value = "17_ROMERO"
$ws1.Range("D479").Value = 115
$ws1.Range("C480").Value = "11_ETCHEVERRY"
$ws1.Range("B481").Value = "22:08"
$ws1.Range("C481").Value = "17_ROMERO"
$ws1.Range("D481").Value = 94
$ws1.Range("A482").Value = "20:47:24"
$ws1.Range("B482").Value = "22:15"
$ws1.Range("C482").Value = "26_HERNANDEZ"
$ws1.Range("D482").Value = 88
$ws1.Range("B483").Value = "22:16"
$ws1.Range("C483").Value = "26_HERNANDEZ"
$ws1.Range("D483").Value = 102
$ws1.Range("B484").Value = "22:21"
$ws1.Range("C484").Value = "23_HERNANDEZ"
$ws1.Range("D484").Value = 107
$ws1.Range("A485").Value = "20:47:24"
$ws1.Range("B485").Value = "22:25"
$ws1.Range("C485").Value = "23_HERNANDEZ"
$ws1.Range("D485").Value = 98
$ws1.Range("E485").Value = "LP1912"
$ws1.Range("A486").Value = "20:34:16"
$ws1.Range("B486").Value = "22:28"
$ws1.Range("C486").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D486").Value = 114
$ws1.Range("E486").Value = "LP1912"
$ws1.Range("A487").Value = "20:34:16"
$ws1.Range("B487").Value = "22:32"
$ws1.Range("C487").Value = "10_OLMOS"
$ws1.Range("D487").Value = 118
$ws1.Range("E487").Value = "LP1912"
$ws1.Range("A488").Value = "20:47:24"
$ws1.Range("B488").Value = "22:34"
$ws1.Range("C488").Value = "10_OLMOS"
$ws1.Range("D488").Value = 107
$ws1.Range("E488").Value = "LP1912"
$ws1.Range("A489").Value = "20:47:24"
$ws1.Range("B489").Value = "22:39"
$ws1.Range("C489").Value = "215A_EL PATO"
$ws1.Range("D489").Value = 112
$ws1.Range("E489").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2.Range("A2").Value = "Última actualización: 20:47:24"
$ws2.Range("A3").Value = "Total filas: 47"
$ws2.Range("A52").Value = "20:47:24"
$ws2.Range("B52").Value = "22:39"
$ws2.Range("C52").Value = "215A_EL PATO"
$ws2.Range("D52").Value = 112
$ws2.Range("E52").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3.Range("A2").Value = "Última actualización: 20:47:24"
$ws3.Range("A3").Value = "Total filas: 62"
$ws3.Range("A66").Value = "20:47:24"
$ws3.Range("B66").Value = "22:20"
$ws3.Range("D66").Value = 93
$ws3.Range("A67").Value = "20:34:16"
$ws3.Range("B67").Value = "22:21"
$ws3.Range("C67").Value = "215B_LP-P MOR-40 Y 115"
$ws3.Range("D67").Value = 107
$ws3.Range("E67").Value = "L6173"

